$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 266-330 (3 new price records were inserted after row 265,
#     shifting the existing records for rows 268-330 down by three rows) ---
$ws.Range("D266").Value = 44476
$ws.Range("K266").Value = "Packham's Triumph"
$ws.Range("L266").Value = "Especial"
$ws.Range("M266").Value = 20
$ws.Range("N266").Value = 275000
$ws.Range("O266").Value = 280000
$ws.Range("P266").Value = 277500
$ws.Range("R266").Value = "Región de O'Higgins"
$ws.Range("S266").Value = 617

$ws.Range("D267").Value = 44476
$ws.Range("K267").Value = "Packham's Triumph"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 20
$ws.Range("N267").Value = 255000
$ws.Range("O267").Value = 260000
$ws.Range("P267").Value = 257500
$ws.Range("R267").Value = "Región de O'Higgins"
$ws.Range("S267").Value = 572

$ws.Range("D268").Value = 44476
$ws.Range("K268").Value = "Packham's Triumph"
$ws.Range("L268").Value = "Segunda"
$ws.Range("M268").Value = 20
$ws.Range("N268").Value = 225000
$ws.Range("O268").Value = 230000
$ws.Range("P268").Value = 227500
$ws.Range("R268").Value = "Región de O'Higgins"
$ws.Range("S268").Value = 506

$ws.Range("D269").Value = 44386
$ws.Range("K269").Value = "Packham's Triumph"
$ws.Range("L269").Value = "Especial"
$ws.Range("M269").Value = 16
$ws.Range("N269").Value = 235000
$ws.Range("O269").Value = 240000
$ws.Range("P269").Value = 237500
$ws.Range("R269").Value = "Región de O'Higgins"
$ws.Range("S269").Value = 528

$ws.Range("D270").Value = 44386
$ws.Range("K270").Value = "Packham's Triumph"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 20
$ws.Range("N270").Value = 205000
$ws.Range("O270").Value = 210000
$ws.Range("P270").Value = 207500
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 461

$ws.Range("D271").Value = 44263
$ws.Range("K271").Value = "Packham's Triumph"
$ws.Range("L271").Value = "Primera"
$ws.Range("M271").Value = 20
$ws.Range("N271").Value = 245000
$ws.Range("O271").Value = 250000
$ws.Range("P271").Value = 247500
$ws.Range("R271").Value = "Región de O'Higgins"
$ws.Range("S271").Value = 550

$ws.Range("D272").Value = 44263
$ws.Range("K272").Value = "Packham's Triumph"
$ws.Range("L272").Value = "Segunda"
$ws.Range("M272").Value = 20
$ws.Range("N272").Value = 205000
$ws.Range("O272").Value = 210000
$ws.Range("P272").Value = 207500
$ws.Range("R272").Value = "Región de O'Higgins"
$ws.Range("S272").Value = 461

$ws.Range("D273").Value = 44306
$ws.Range("K273").Value = "Packham's Triumph"
$ws.Range("L273").Value = "Especial"
$ws.Range("M273").Value = 26
$ws.Range("N273").Value = 255000
$ws.Range("O273").Value = 260000
$ws.Range("P273").Value = 257500
$ws.Range("R273").Value = "Región de O'Higgins"
$ws.Range("S273").Value = 572

$ws.Range("D274").Value = 44306
$ws.Range("K274").Value = "Packham's Triumph"
$ws.Range("L274").Value = "Primera"
$ws.Range("M274").Value = 20
$ws.Range("N274").Value = 235000
$ws.Range("O274").Value = 240000
$ws.Range("P274").Value = 237500
$ws.Range("R274").Value = "Región de O'Higgins"
$ws.Range("S274").Value = 528

$ws.Range("D275").Value = 44369
$ws.Range("K275").Value = "Packham's Triumph"
$ws.Range("L275").Value = "Especial"
$ws.Range("M275").Value = 24
$ws.Range("N275").Value = 225000
$ws.Range("O275").Value = 230000
$ws.Range("P275").Value = 227500
$ws.Range("R275").Value = "Región de O'Higgins"
$ws.Range("S275").Value = 506

$ws.Range("D276").Value = 44369
$ws.Range("K276").Value = "Packham's Triumph"
$ws.Range("L276").Value = "Primera"
$ws.Range("M276").Value = 20
$ws.Range("N276").Value = 205000
$ws.Range("O276").Value = 210000
$ws.Range("P276").Value = 207500
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value = 461

$ws.Range("D277").Value = 44369
$ws.Range("K277").Value = "Packham's Triumph"
$ws.Range("L277").Value = "Segunda"
$ws.Range("M277").Value = 18
$ws.Range("N277").Value = 185000
$ws.Range("O277").Value = 190000
$ws.Range("P277").Value = 187500
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 417

$ws.Range("D278").Value = 44301
$ws.Range("K278").Value = "Packham's Triumph"
$ws.Range("L278").Value = "Primera"
$ws.Range("M278").Value = 20
$ws.Range("N278").Value = 225000
$ws.Range("O278").Value = 230000
$ws.Range("P278").Value = 227500
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 506

$ws.Range("D279").Value = 44301
$ws.Range("K279").Value = "Packham's Triumph"
$ws.Range("L279").Value = "Segunda"
$ws.Range("M279").Value = 20
$ws.Range("N279").Value = 205000
$ws.Range("O279").Value = 210000
$ws.Range("P279").Value = 207500
$ws.Range("R279").Value = "Región de O'Higgins"
$ws.Range("S279").Value = 461

$ws.Range("D280").Value = 44301
$ws.Range("K280").Value = "Winter Nelis"
$ws.Range("L280").Value = "Primera"
$ws.Range("M280").Value = 20
$ws.Range("N280").Value = 235000
$ws.Range("O280").Value = 240000
$ws.Range("P280").Value = 237500
$ws.Range("R280").Value = "Región de O'Higgins"
$ws.Range("S280").Value = 528

$ws.Range("D281").Value = 44301
$ws.Range("K281").Value = "Winter Nelis"
$ws.Range("L281").Value = "Segunda"
$ws.Range("M281").Value = 20
$ws.Range("N281").Value = 215000
$ws.Range("O281").Value = 220000
$ws.Range("P281").Value = 217500
$ws.Range("R281").Value = "Región de O'Higgins"
$ws.Range("S281").Value = 483

$ws.Range("D282").Value = 44357
$ws.Range("K282").Value = "Packham's Triumph"
$ws.Range("L282").Value = "Especial"
$ws.Range("M282").Value = 16
$ws.Range("N282").Value = 245000
$ws.Range("O282").Value = 250000
$ws.Range("P282").Value = 247500
$ws.Range("R282").Value = "Región de O'Higgins"
$ws.Range("S282").Value = 550

$ws.Range("D283").Value = 44357
$ws.Range("K283").Value = "Packham's Triumph"
$ws.Range("L283").Value = "Primera"
$ws.Range("M283").Value = 20
$ws.Range("N283").Value = 215000
$ws.Range("O283").Value = 220000
$ws.Range("P283").Value = 217500
$ws.Range("R283").Value = "Región de O'Higgins"
$ws.Range("S283").Value = 483

$ws.Range("D284").Value = 44357
$ws.Range("K284").Value = "Packham's Triumph"
$ws.Range("L284").Value = "Segunda"
$ws.Range("M284").Value = 20
$ws.Range("N284").Value = 185000
$ws.Range("O284").Value = 190000
$ws.Range("P284").Value = 187500
$ws.Range("R284").Value = "Región de O'Higgins"
$ws.Range("S284").Value = 417

$ws.Range("D285").Value = 44357
$ws.Range("K285").Value = "Winter Nelis"
$ws.Range("L285").Value = "Primera"
$ws.Range("M285").Value = 20
$ws.Range("N285").Value = 190000
$ws.Range("O285").Value = 195000
$ws.Range("P285").Value = 192500
$ws.Range("R285").Value = "Región de O'Higgins"
$ws.Range("S285").Value = 428

$ws.Range("D286").Value = 44357
$ws.Range("K286").Value = "Winter Nelis"
$ws.Range("L286").Value = "Segunda"
$ws.Range("M286").Value = 20
$ws.Range("N286").Value = 170000
$ws.Range("O286").Value = 175000
$ws.Range("P286").Value = 172500
$ws.Range("R286").Value = "Región de O'Higgins"
$ws.Range("S286").Value = 383

$ws.Range("D287").Value = 44298
$ws.Range("K287").Value = "Packham's Triumph"
$ws.Range("L287").Value = "Primera"
$ws.Range("M287").Value = 20
$ws.Range("N287").Value = 225000
$ws.Range("O287").Value = 230000
$ws.Range("P287").Value = 227500
$ws.Range("R287").Value = "Región de O'Higgins"
$ws.Range("S287").Value = 506

$ws.Range("D288").Value = 44298
$ws.Range("K288").Value = "Packham's Triumph"
$ws.Range("L288").Value = "Segunda"
$ws.Range("M288").Value = 20
$ws.Range("N288").Value = 205000
$ws.Range("O288").Value = 210000
$ws.Range("P288").Value = 207500
$ws.Range("R288").Value = "Región de O'Higgins"
$ws.Range("S288").Value = 461

$ws.Range("D289").Value = 44397
$ws.Range("K289").Value = "Packham's Triumph"
$ws.Range("L289").Value = "Especial"
$ws.Range("M289").Value = 24
$ws.Range("N289").Value = 220000
$ws.Range("O289").Value = 230000
$ws.Range("P289").Value = 225000
$ws.Range("R289").Value = "Provincia de Curicó"
$ws.Range("S289").Value = 500

$ws.Range("D290").Value = 44397
$ws.Range("K290").Value = "Packham's Triumph"
$ws.Range("L290").Value = "Primera"
$ws.Range("M290").Value = 20
$ws.Range("N290").Value = 200000
$ws.Range("O290").Value = 210000
$ws.Range("P290").Value = 205000
$ws.Range("R290").Value = "Provincia de Curicó"
$ws.Range("S290").Value = 456

$ws.Range("D291").Value = 44397
$ws.Range("K291").Value = "Packham's Triumph"
$ws.Range("L291").Value = "Segunda"
$ws.Range("M291").Value = 18
$ws.Range("N291").Value = 180000
$ws.Range("O291").Value = 190000
$ws.Range("P291").Value = 185000
$ws.Range("R291").Value = "Provincia de Curicó"
$ws.Range("S291").Value = 411

$ws.Range("D292").Value = 44397
$ws.Range("K292").Value = "Winter Nelis"
$ws.Range("L292").Value = "Especial"
$ws.Range("M292").Value = 20
$ws.Range("N292").Value = 230000
$ws.Range("O292").Value = 240000
$ws.Range("P292").Value = 235000
$ws.Range("R292").Value = "Provincia de Curicó"
$ws.Range("S292").Value = 522

$ws.Range("D293").Value = 44397
$ws.Range("K293").Value = "Winter Nelis"
$ws.Range("L293").Value = "Primera"
$ws.Range("M293").Value = 16
$ws.Range("N293").Value = 210000
$ws.Range("O293").Value = 220000
$ws.Range("P293").Value = 215000
$ws.Range("R293").Value = "Provincia de Curicó"
$ws.Range("S293").Value = 478

$ws.Range("D294").Value = 44397
$ws.Range("K294").Value = "Winter Nelis"
$ws.Range("L294").Value = "Segunda"
$ws.Range("M294").Value = 14
$ws.Range("N294").Value = 170000
$ws.Range("O294").Value = 180000
$ws.Range("P294").Value = 175000
$ws.Range("R294").Value = "Provincia de Curicó"
$ws.Range("S294").Value = 389

$ws.Range("D295").Value = 44414
$ws.Range("K295").Value = "Packham's Triumph"
$ws.Range("L295").Value = "Especial"
$ws.Range("M295").Value = 16
$ws.Range("N295").Value = 260000
$ws.Range("O295").Value = 265000
$ws.Range("P295").Value = 262500
$ws.Range("R295").Value = "Región de O'Higgins"
$ws.Range("S295").Value = 583

$ws.Range("D296").Value = 44414
$ws.Range("K296").Value = "Packham's Triumph"
$ws.Range("L296").Value = "Primera"
$ws.Range("M296").Value = 16
$ws.Range("N296").Value = 230000
$ws.Range("O296").Value = 235000
$ws.Range("P296").Value = 232500
$ws.Range("R296").Value = "Región de O'Higgins"
$ws.Range("S296").Value = 517

$ws.Range("D297").Value = 44414
$ws.Range("K297").Value = "Packham's Triumph"
$ws.Range("L297").Value = "Segunda"
$ws.Range("M297").Value = 16
$ws.Range("N297").Value = 210000
$ws.Range("O297").Value = 215000
$ws.Range("P297").Value = 212500
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 472

$ws.Range("D298").Value = 44414
$ws.Range("K298").Value = "Winter Nelis"
$ws.Range("L298").Value = "Especial"
$ws.Range("M298").Value = 10
$ws.Range("N298").Value = 235000
$ws.Range("O298").Value = 240000
$ws.Range("P298").Value = 237500
$ws.Range("R298").Value = "Provincia de Curicó"
$ws.Range("S298").Value = 528

$ws.Range("D299").Value = 44414
$ws.Range("K299").Value = "Winter Nelis"
$ws.Range("L299").Value = "Primera"
$ws.Range("M299").Value = 20
$ws.Range("N299").Value = 215000
$ws.Range("O299").Value = 220000
$ws.Range("P299").Value = 217500
$ws.Range("R299").Value = "Provincia de Curicó"
$ws.Range("S299").Value = 483

$ws.Range("D300").Value = 44414
$ws.Range("K300").Value = "Winter Nelis"
$ws.Range("L300").Value = "Segunda"
$ws.Range("M300").Value = 10
$ws.Range("N300").Value = 195000
$ws.Range("O300").Value = 200000
$ws.Range("P300").Value = 197500
$ws.Range("R300").Value = "Provincia de Curicó"
$ws.Range("S300").Value = 439

$ws.Range("D301").Value = 44333
$ws.Range("K301").Value = "Packham's Triumph"
$ws.Range("L301").Value = "Especial"
$ws.Range("M301").Value = 20
$ws.Range("N301").Value = 225000
$ws.Range("O301").Value = 230000
$ws.Range("P301").Value = 227500
$ws.Range("R301").Value = "Región de O'Higgins"
$ws.Range("S301").Value = 506

$ws.Range("D302").Value = 44333
$ws.Range("K302").Value = "Packham's Triumph"
$ws.Range("L302").Value = "Primera"
$ws.Range("M302").Value = 20
$ws.Range("N302").Value = 205000
$ws.Range("O302").Value = 210000
$ws.Range("P302").Value = 207500
$ws.Range("R302").Value = "Región de O'Higgins"
$ws.Range("S302").Value = 461

$ws.Range("D303").Value = 44333
$ws.Range("K303").Value = "Packham's Triumph"
$ws.Range("L303").Value = "Segunda"
$ws.Range("M303").Value = 20
$ws.Range("N303").Value = 185000
$ws.Range("O303").Value = 190000
$ws.Range("P303").Value = 187500
$ws.Range("R303").Value = "Región de O'Higgins"
$ws.Range("S303").Value = 417

$ws.Range("D304").Value = 44392
$ws.Range("K304").Value = "Packham's Triumph"
$ws.Range("L304").Value = "Especial"
$ws.Range("M304").Value = 20
$ws.Range("N304").Value = 235000
$ws.Range("O304").Value = 240000
$ws.Range("P304").Value = 237500
$ws.Range("R304").Value = "Región de O'Higgins"
$ws.Range("S304").Value = 528

$ws.Range("D305").Value = 44392
$ws.Range("K305").Value = "Packham's Triumph"
$ws.Range("L305").Value = "Primera"
$ws.Range("M305").Value = 20
$ws.Range("N305").Value = 205000
$ws.Range("O305").Value = 210000
$ws.Range("P305").Value = 207500
$ws.Range("R305").Value = "Región de O'Higgins"
$ws.Range("S305").Value = 461

$ws.Range("D306").Value = 44392
$ws.Range("K306").Value = "Packham's Triumph"
$ws.Range("L306").Value = "Segunda"
$ws.Range("M306").Value = 20
$ws.Range("N306").Value = 185000
$ws.Range("O306").Value = 190000
$ws.Range("P306").Value = 187500
$ws.Range("R306").Value = "Región de O'Higgins"
$ws.Range("S306").Value = 417

$ws.Range("D307").Value = 44392
$ws.Range("K307").Value = "Winter Nelis"
$ws.Range("L307").Value = "Especial"
$ws.Range("M307").Value = 16
$ws.Range("N307").Value = 225000
$ws.Range("O307").Value = 230000
$ws.Range("P307").Value = 227500
$ws.Range("R307").Value = "Región de O'Higgins"
$ws.Range("S307").Value = 506

$ws.Range("D308").Value = 44392
$ws.Range("K308").Value = "Winter Nelis"
$ws.Range("L308").Value = "Primera"
$ws.Range("M308").Value = 20
$ws.Range("N308").Value = 205000
$ws.Range("O308").Value = 210000
$ws.Range("P308").Value = 207500
$ws.Range("R308").Value = "Región de O'Higgins"
$ws.Range("S308").Value = 461

$ws.Range("D309").Value = 44392
$ws.Range("K309").Value = "Winter Nelis"
$ws.Range("L309").Value = "Segunda"
$ws.Range("M309").Value = 16
$ws.Range("N309").Value = 185000
$ws.Range("O309").Value = 190000
$ws.Range("P309").Value = 187500
$ws.Range("R309").Value = "Región de O'Higgins"
$ws.Range("S309").Value = 417

$ws.Range("D310").Value = 44425
$ws.Range("K310").Value = "Packham's Triumph"
$ws.Range("L310").Value = "Especial"
$ws.Range("M310").Value = 26
$ws.Range("N310").Value = 250000
$ws.Range("O310").Value = 260000
$ws.Range("P310").Value = 255000
$ws.Range("R310").Value = "Provincia de Curicó"
$ws.Range("S310").Value = 567

$ws.Range("D311").Value = 44425
$ws.Range("K311").Value = "Packham's Triumph"
$ws.Range("L311").Value = "Primera"
$ws.Range("M311").Value = 20
$ws.Range("N311").Value = 230000
$ws.Range("O311").Value = 240000
$ws.Range("P311").Value = 235000
$ws.Range("R311").Value = "Provincia de Curicó"
$ws.Range("S311").Value = 522

$ws.Range("D312").Value = 44425
$ws.Range("K312").Value = "Packham's Triumph"
$ws.Range("L312").Value = "Segunda"
$ws.Range("M312").Value = 18
$ws.Range("N312").Value = 200000
$ws.Range("O312").Value = 210000
$ws.Range("P312").Value = 205000
$ws.Range("R312").Value = "Provincia de Curicó"
$ws.Range("S312").Value = 456

$ws.Range("D313").Value = 44425
$ws.Range("K313").Value = "Winter Nelis"
$ws.Range("L313").Value = "Primera"
$ws.Range("M313").Value = 24
$ws.Range("N313").Value = 230000
$ws.Range("O313").Value = 240000
$ws.Range("P313").Value = 235000
$ws.Range("R313").Value = "Provincia de Curicó"
$ws.Range("S313").Value = 522

$ws.Range("D314").Value = 44425
$ws.Range("K314").Value = "Winter Nelis"
$ws.Range("L314").Value = "Segunda"
$ws.Range("M314").Value = 18
$ws.Range("N314").Value = 210000
$ws.Range("O314").Value = 220000
$ws.Range("P314").Value = 215000
$ws.Range("R314").Value = "Provincia de Curicó"
$ws.Range("S314").Value = 478

$ws.Range("D315").Value = 44390
$ws.Range("K315").Value = "Packham's Triumph"
$ws.Range("L315").Value = "Especial"
$ws.Range("M315").Value = 24
$ws.Range("N315").Value = 235000
$ws.Range("O315").Value = 240000
$ws.Range("P315").Value = 237500
$ws.Range("R315").Value = "Región de O'Higgins"
$ws.Range("S315").Value = 528

$ws.Range("D316").Value = 44390
$ws.Range("K316").Value = "Packham's Triumph"
$ws.Range("L316").Value = "Primera"
$ws.Range("M316").Value = 20
$ws.Range("N316").Value = 205000
$ws.Range("O316").Value = 210000
$ws.Range("P316").Value = 207500
$ws.Range("R316").Value = "Región de O'Higgins"
$ws.Range("S316").Value = 461

$ws.Range("D317").Value = 44390
$ws.Range("K317").Value = "Packham's Triumph"
$ws.Range("L317").Value = "Segunda"
$ws.Range("M317").Value = 18
$ws.Range("N317").Value = 185000
$ws.Range("O317").Value = 190000
$ws.Range("P317").Value = 187500
$ws.Range("R317").Value = "Región de O'Higgins"
$ws.Range("S317").Value = 417

$ws.Range("D318").Value = 44466
$ws.Range("K318").Value = "Packham's Triumph"
$ws.Range("L318").Value = "Especial"
$ws.Range("M318").Value = 16
$ws.Range("N318").Value = 285000
$ws.Range("O318").Value = 290000
$ws.Range("P318").Value = 287500
$ws.Range("R318").Value = "Región de O'Higgins"
$ws.Range("S318").Value = 639

$ws.Range("D319").Value = 44466
$ws.Range("K319").Value = "Packham's Triumph"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 16
$ws.Range("N319").Value = 255000
$ws.Range("O319").Value = 260000
$ws.Range("P319").Value = 257500
$ws.Range("R319").Value = "Región de O'Higgins"
$ws.Range("S319").Value = 572

$ws.Range("D320").Value = 44466
$ws.Range("K320").Value = "Packham's Triumph"
$ws.Range("L320").Value = "Segunda"
$ws.Range("M320").Value = 20
$ws.Range("N320").Value = 235000
$ws.Range("O320").Value = 240000
$ws.Range("P320").Value = 237500
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 528

$ws.Range("D321").Value = 44438
$ws.Range("K321").Value = "Packham's Triumph"
$ws.Range("L321").Value = "Especial"
$ws.Range("M321").Value = 16
$ws.Range("N321").Value = 270000
$ws.Range("O321").Value = 280000
$ws.Range("P321").Value = 275000
$ws.Range("R321").Value = "Región de O'Higgins"
$ws.Range("S321").Value = 611

$ws.Range("D322").Value = 44438
$ws.Range("K322").Value = "Packham's Triumph"
$ws.Range("L322").Value = "Primera"
$ws.Range("M322").Value = 20
$ws.Range("N322").Value = 250000
$ws.Range("O322").Value = 260000
$ws.Range("P322").Value = 255000
$ws.Range("R322").Value = "Región de O'Higgins"
$ws.Range("S322").Value = 567

$ws.Range("D323").Value = 44438
$ws.Range("K323").Value = "Packham's Triumph"
$ws.Range("L323").Value = "Segunda"
$ws.Range("M323").Value = 20
$ws.Range("N323").Value = 220000
$ws.Range("O323").Value = 230000
$ws.Range("P323").Value = 225000
$ws.Range("R323").Value = "Región de O'Higgins"
$ws.Range("S323").Value = 500

$ws.Range("D324").Value = 44286
$ws.Range("K324").Value = "Winter Nelis"
$ws.Range("L324").Value = "Primera"
$ws.Range("M324").Value = 22
$ws.Range("N324").Value = 195000
$ws.Range("O324").Value = 200000
$ws.Range("P324").Value = 197500
$ws.Range("R324").Value = "Región de O'Higgins"
$ws.Range("S324").Value = 439

$ws.Range("D325").Value = 44286
$ws.Range("K325").Value = "Winter Nelis"
$ws.Range("L325").Value = "Segunda"
$ws.Range("M325").Value = 18
$ws.Range("N325").Value = 155000
$ws.Range("O325").Value = 160000
$ws.Range("P325").Value = 157500
$ws.Range("R325").Value = "Región de O'Higgins"
$ws.Range("S325").Value = 350

$ws.Range("D326").Value = 44389
$ws.Range("K326").Value = "Packham's Triumph"
$ws.Range("L326").Value = "Especial"
$ws.Range("M326").Value = 16
$ws.Range("N326").Value = 235000
$ws.Range("O326").Value = 240000
$ws.Range("P326").Value = 237500
$ws.Range("R326").Value = "Región de O'Higgins"
$ws.Range("S326").Value = 528

$ws.Range("D327").Value = 44389
$ws.Range("K327").Value = "Packham's Triumph"
$ws.Range("L327").Value = "Primera"
$ws.Range("M327").Value = 20
$ws.Range("N327").Value = 205000
$ws.Range("O327").Value = 210000
$ws.Range("P327").Value = 207500
$ws.Range("R327").Value = "Región de O'Higgins"
$ws.Range("S327").Value = 461

$ws.Range("D328").Value = 44389
$ws.Range("K328").Value = "Packham's Triumph"
$ws.Range("L328").Value = "Segunda"
$ws.Range("M328").Value = 16
$ws.Range("N328").Value = 185000
$ws.Range("O328").Value = 190000
$ws.Range("P328").Value = 187500
$ws.Range("R328").Value = "Región de O'Higgins"
$ws.Range("S328").Value = 417

$ws.Range("D329").Value = 44312
$ws.Range("K329").Value = "Packham's Triumph"
$ws.Range("L329").Value = "Especial"
$ws.Range("M329").Value = 20
$ws.Range("N329").Value = 255000
$ws.Range("O329").Value = 260000
$ws.Range("P329").Value = 257500
$ws.Range("R329").Value = "Región de O'Higgins"
$ws.Range("S329").Value = 572

$ws.Range("D330").Value = 44312
$ws.Range("K330").Value = "Packham's Triumph"
$ws.Range("L330").Value = "Primera"
$ws.Range("M330").Value = 20
$ws.Range("N330").Value = 235000
$ws.Range("O330").Value = 240000
$ws.Range("P330").Value = 237500
$ws.Range("R330").Value = "Región de O'Higgins"
$ws.Range("S330").Value = 528

# --- Add new rows 331-333 (copy formatting/values from row 330, one row at a time, then set the real values) ---
$ws.Range("A330:T330").Copy($ws.Range("A331:T331"))
$ws.Range("A330:T330").Copy($ws.Range("A332:T332"))
$ws.Range("A330:T330").Copy($ws.Range("A333:T333"))

$ws.Range("D331").Value = 44312
$ws.Range("K331").Value = "Packham's Triumph"
$ws.Range("L331").Value = "Segunda"
$ws.Range("M331").Value = 20
$ws.Range("N331").Value = 195000
$ws.Range("O331").Value = 200000
$ws.Range("P331").Value = 197500
$ws.Range("R331").Value = "Región de O'Higgins"
$ws.Range("S331").Value = 439

$ws.Range("D332").Value = 44326
$ws.Range("K332").Value = "Packham's Triumph"
$ws.Range("L332").Value = "Primera"
$ws.Range("M332").Value = 20
$ws.Range("N332").Value = 215000
$ws.Range("O332").Value = 220000
$ws.Range("P332").Value = 217500
$ws.Range("R332").Value = "Región de O'Higgins"
$ws.Range("S332").Value = 483

$ws.Range("D333").Value = 44326
$ws.Range("K333").Value = "Packham's Triumph"
$ws.Range("L333").Value = "Segunda"
$ws.Range("M333").Value = 20
$ws.Range("N333").Value = 185000
$ws.Range("O333").Value = 190000
$ws.Range("P333").Value = 187500
$ws.Range("R333").Value = "Región de O'Higgins"
$ws.Range("S333").Value = 417
